$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------
# Helper: write a plain text value into a cell without Excel's
# auto-conversion of date-looking strings ("2026-02-03") into date
# serial numbers. We stage the text (forced via a text number format)
# in a scratch cell, then bring only the VALUE over with paste-special,
# leaving whatever formatting is already on the destination untouched.
# ----------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

# 1) Fill in the previously blank clock-out time & duration on row 15
$ws.Range("C15").Value = "18:49:35"
$ws.Range("D15").Value = "2.57 Hours"

# 2) Add new row 16 (full entry) - copy formatting from row 15 first
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

Set-TextValue $ws.Range("A16") "2026-02-03"
$ws.Range("B16").Value = "20:26:35"
$ws.Range("C16").Value = "21:25:16"
$ws.Range("D16").Value = "0.98 Hours"

# 3) Add new row 17 (clock-in only; clock-out/duration remain blank)
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

Set-TextValue $ws.Range("A17") "2026-02-03"
$ws.Range("B17").Value = "21:32:50"

$excel.CutCopyMode = 0
